# Fruta / hortaliza, semanal
# This edit re-shuffles the weekly records: for each data row (2..46) the
# values in columns D (Fecha), M (Volumen), N (Precio minimo),
# O (Precio maximo), P (Precio promedio ponderado) and S (Precio $/Kg)
# are replaced by the values that used to live in a different row.
# All the other columns (A,B,C,E,F,G,H,I,J,K,L,Q,R,T) are constant across
# the whole table and are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# destination row -> source row (source row's D/M/N/O/P/S values are moved
# into destination row)
$map = @{
    2 = 44
    3 = 36
    4 = 32
    5 = 24
    6 = 46
    7 = 7
    8 = 45
    9 = 29
    10 = 38
    11 = 34
    12 = 40
    13 = 21
    14 = 6
    15 = 4
    16 = 41
    17 = 15
    18 = 10
    19 = 42
    20 = 13
    21 = 12
    22 = 30
    23 = 9
    24 = 3
    25 = 22
    26 = 11
    27 = 2
    28 = 8
    29 = 33
    30 = 39
    31 = 14
    32 = 16
    33 = 35
    34 = 43
    35 = 27
    36 = 20
    37 = 37
    38 = 28
    39 = 26
    40 = 31
    41 = 5
    42 = 18
    43 = 23
    44 = 17
    45 = 25
    46 = 19
}

$cols = @("D", "M", "N", "O", "P", "S")

# Snapshot the current value of every relevant cell before writing anything,
# since several rows swap values with each other.
$snapshot = @{}
for ($r = 2; $r -le 46; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $rowVals
}

foreach ($dstRow in $map.Keys) {
    $srcRow = $map[$dstRow]
    $srcVals = $snapshot[$srcRow]
    foreach ($c in $cols) {
        $ws.Range("$c$dstRow").Value = $srcVals[$c]
    }
}
